$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = '281474990207675-1749682514529'
$ws.Range("B2").Value = 'Defensive Driving'
$ws.Range("C2").Value = '2025-06-11T16:55:14.529'
$ws.Range("D2").Value = '''281474990207675'
$ws.Range("E2").Value = '''122'
$ws.Range("F2").Value = '''51833996'
$ws.Range("G2").Value = 'MARTÍN QUEZADA'
$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990207675/1749682509529/zqlwK29jaz-camera-video-segment-driver-1749682514529.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPKB53JIC%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T170058Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBcaCXVzLXdlc3QtMiJIMEYCIQD6QllGW0%2Bq0TVWZMTrFyq%2FzctklkBpmTevR5cigwbYpAIhAM%2BaFQ4HkAWSYncq2J05pvuQjHQugDSlPakTP9CUrk9pKuYDCPD%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwtArUTuLNw%2BvXf9b0qugMkH9CxTIBgei6Y1mVM1SKBEL%2FYxLhENIC5xiLffQRvKfcuOuULFPgrg4zpl6BW%2FSEW938AsTsnzcZyb%2F5F7lOCfEgPkM0XlBiLU6mPQhf1m96y9ESQ7hWXAtbVPReGsK5932UXmxHdAt%2BrvGeRb0zs7uD1p5ipLNrO8ep%2FFTI0VQP%2F7ll1Mia7TyvGLSOVR6h9K8V9udrWgToYtYWmLM8h3GM%2FFfc95vlnqU9zIAEO4thgsRs2WsBLSaOH8CTsD3g%2BoyzM0Xe2ogPMYLf2MJPknQrpTWa%2BRChMUVM8JyJhOBlimM7HnCHn84ty9A6FR63ndE15ze7Eyj%2F8Cw3xY8%2FTpP536z0UyF8QF6fMpc9sxTM%2BrlLbv231Cx1ngj0PckQpSB7ssg3mzUfQmwXTMpx0jsvflSfIA%2BaxGCbb7WoijhNuPDFcdGzPiF5%2FQd6e8mkjs5ShTZA1ZgsKFTCVZclxgZcHnPfNwu0nlCY%2BF01TP7y4Kx4ce456zAFsJ95rPHtgQLIekFNhNuJbNUGxtMdILc8FGSeKzJUScJrR%2BHh1Jmmr0GAJmNTXtM3jF6ggEgeU0nPuoXoovgqIMOTUq8IGOqQBEb8nZUZBCOuV6XQ3C%2BQXTY13rYR2anHJIiQJIP49O4VU81pKrgr0SqCkpntZgHaQpTowDdsjqZoBjsQNsMwFZ4wHTWXSrzoHZzKulSj5OgHruJX910c0ejMRikZvW7jYncLaD3pgxEO9ssLjhaIQG8RA%2FCPwovT8h00ZIZ0ZcHki%2F5rs4c5BJbMinhyN%2FwqtHr%2BAUpPZWkTHyLw0pkQm7sDS%2BR0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2001%3A00%3A58%20GMT&X-Amz-Signature=6d39687af7441579be9eaf9069d9c47dbc845afbfb07a4cc8984a1f5992889bc'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990207675/1749682509529/JfEAwSUnRy-camera-video-segment-1749682514529.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPKB53JIC%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T170058Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBcaCXVzLXdlc3QtMiJIMEYCIQD6QllGW0%2Bq0TVWZMTrFyq%2FzctklkBpmTevR5cigwbYpAIhAM%2BaFQ4HkAWSYncq2J05pvuQjHQugDSlPakTP9CUrk9pKuYDCPD%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwtArUTuLNw%2BvXf9b0qugMkH9CxTIBgei6Y1mVM1SKBEL%2FYxLhENIC5xiLffQRvKfcuOuULFPgrg4zpl6BW%2FSEW938AsTsnzcZyb%2F5F7lOCfEgPkM0XlBiLU6mPQhf1m96y9ESQ7hWXAtbVPReGsK5932UXmxHdAt%2BrvGeRb0zs7uD1p5ipLNrO8ep%2FFTI0VQP%2F7ll1Mia7TyvGLSOVR6h9K8V9udrWgToYtYWmLM8h3GM%2FFfc95vlnqU9zIAEO4thgsRs2WsBLSaOH8CTsD3g%2BoyzM0Xe2ogPMYLf2MJPknQrpTWa%2BRChMUVM8JyJhOBlimM7HnCHn84ty9A6FR63ndE15ze7Eyj%2F8Cw3xY8%2FTpP536z0UyF8QF6fMpc9sxTM%2BrlLbv231Cx1ngj0PckQpSB7ssg3mzUfQmwXTMpx0jsvflSfIA%2BaxGCbb7WoijhNuPDFcdGzPiF5%2FQd6e8mkjs5ShTZA1ZgsKFTCVZclxgZcHnPfNwu0nlCY%2BF01TP7y4Kx4ce456zAFsJ95rPHtgQLIekFNhNuJbNUGxtMdILc8FGSeKzJUScJrR%2BHh1Jmmr0GAJmNTXtM3jF6ggEgeU0nPuoXoovgqIMOTUq8IGOqQBEb8nZUZBCOuV6XQ3C%2BQXTY13rYR2anHJIiQJIP49O4VU81pKrgr0SqCkpntZgHaQpTowDdsjqZoBjsQNsMwFZ4wHTWXSrzoHZzKulSj5OgHruJX910c0ejMRikZvW7jYncLaD3pgxEO9ssLjhaIQG8RA%2FCPwovT8h00ZIZ0ZcHki%2F5rs4c5BJbMinhyN%2FwqtHr%2BAUpPZWkTHyLw0pkQm7sDS%2BR0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2001%3A00%3A58%20GMT&X-Amz-Signature=39282a49e069c939ff47c9eb2bc54cfe6748f5e5c89ea716e267044ebc1091cb'
$ws.Range("H2").Value = 20.64991883
$ws.Range("I2").Value = -103.36675577
$ws.Range("J2").Value = 0.6070539355278015

# Row 3
$ws.Range("A3").Value = '281474991154589-1749681712753'
$ws.Range("B3").Value = 'Defensive Driving'
$ws.Range("C3").Value = '2025-06-11T16:41:52.753'
$ws.Range("D3").Value = '''281474991154589'
$ws.Range("E3").Value = '''146'
$ws.Range("F3").Value = '''51834110'
$ws.Range("G3").Value = 'ALEJANDRO SUAREZ QUEZADA'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991154589/1749681707753/UtIzP4FPML-camera-video-segment-driver-1749681712753.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPKB53JIC%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T170058Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBcaCXVzLXdlc3QtMiJIMEYCIQD6QllGW0%2Bq0TVWZMTrFyq%2FzctklkBpmTevR5cigwbYpAIhAM%2BaFQ4HkAWSYncq2J05pvuQjHQugDSlPakTP9CUrk9pKuYDCPD%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwtArUTuLNw%2BvXf9b0qugMkH9CxTIBgei6Y1mVM1SKBEL%2FYxLhENIC5xiLffQRvKfcuOuULFPgrg4zpl6BW%2FSEW938AsTsnzcZyb%2F5F7lOCfEgPkM0XlBiLU6mPQhf1m96y9ESQ7hWXAtbVPReGsK5932UXmxHdAt%2BrvGeRb0zs7uD1p5ipLNrO8ep%2FFTI0VQP%2F7ll1Mia7TyvGLSOVR6h9K8V9udrWgToYtYWmLM8h3GM%2FFfc95vlnqU9zIAEO4thgsRs2WsBLSaOH8CTsD3g%2BoyzM0Xe2ogPMYLf2MJPknQrpTWa%2BRChMUVM8JyJhOBlimM7HnCHn84ty9A6FR63ndE15ze7Eyj%2F8Cw3xY8%2FTpP536z0UyF8QF6fMpc9sxTM%2BrlLbv231Cx1ngj0PckQpSB7ssg3mzUfQmwXTMpx0jsvflSfIA%2BaxGCbb7WoijhNuPDFcdGzPiF5%2FQd6e8mkjs5ShTZA1ZgsKFTCVZclxgZcHnPfNwu0nlCY%2BF01TP7y4Kx4ce456zAFsJ95rPHtgQLIekFNhNuJbNUGxtMdILc8FGSeKzJUScJrR%2BHh1Jmmr0GAJmNTXtM3jF6ggEgeU0nPuoXoovgqIMOTUq8IGOqQBEb8nZUZBCOuV6XQ3C%2BQXTY13rYR2anHJIiQJIP49O4VU81pKrgr0SqCkpntZgHaQpTowDdsjqZoBjsQNsMwFZ4wHTWXSrzoHZzKulSj5OgHruJX910c0ejMRikZvW7jYncLaD3pgxEO9ssLjhaIQG8RA%2FCPwovT8h00ZIZ0ZcHki%2F5rs4c5BJbMinhyN%2FwqtHr%2BAUpPZWkTHyLw0pkQm7sDS%2BR0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2001%3A00%3A58%20GMT&X-Amz-Signature=93b2be2b6df50d3b837be5a3281d406eaa42dfb9dec27a06f4469beccfa33e8b'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991154589/1749681707753/418ikFdSWa-camera-video-segment-1749681712753.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSPKB53JIC%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T170058Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBcaCXVzLXdlc3QtMiJIMEYCIQD6QllGW0%2Bq0TVWZMTrFyq%2FzctklkBpmTevR5cigwbYpAIhAM%2BaFQ4HkAWSYncq2J05pvuQjHQugDSlPakTP9CUrk9pKuYDCPD%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgwtArUTuLNw%2BvXf9b0qugMkH9CxTIBgei6Y1mVM1SKBEL%2FYxLhENIC5xiLffQRvKfcuOuULFPgrg4zpl6BW%2FSEW938AsTsnzcZyb%2F5F7lOCfEgPkM0XlBiLU6mPQhf1m96y9ESQ7hWXAtbVPReGsK5932UXmxHdAt%2BrvGeRb0zs7uD1p5ipLNrO8ep%2FFTI0VQP%2F7ll1Mia7TyvGLSOVR6h9K8V9udrWgToYtYWmLM8h3GM%2FFfc95vlnqU9zIAEO4thgsRs2WsBLSaOH8CTsD3g%2BoyzM0Xe2ogPMYLf2MJPknQrpTWa%2BRChMUVM8JyJhOBlimM7HnCHn84ty9A6FR63ndE15ze7Eyj%2F8Cw3xY8%2FTpP536z0UyF8QF6fMpc9sxTM%2BrlLbv231Cx1ngj0PckQpSB7ssg3mzUfQmwXTMpx0jsvflSfIA%2BaxGCbb7WoijhNuPDFcdGzPiF5%2FQd6e8mkjs5ShTZA1ZgsKFTCVZclxgZcHnPfNwu0nlCY%2BF01TP7y4Kx4ce456zAFsJ95rPHtgQLIekFNhNuJbNUGxtMdILc8FGSeKzJUScJrR%2BHh1Jmmr0GAJmNTXtM3jF6ggEgeU0nPuoXoovgqIMOTUq8IGOqQBEb8nZUZBCOuV6XQ3C%2BQXTY13rYR2anHJIiQJIP49O4VU81pKrgr0SqCkpntZgHaQpTowDdsjqZoBjsQNsMwFZ4wHTWXSrzoHZzKulSj5OgHruJX910c0ejMRikZvW7jYncLaD3pgxEO9ssLjhaIQG8RA%2FCPwovT8h00ZIZ0ZcHki%2F5rs4c5BJbMinhyN%2FwqtHr%2BAUpPZWkTHyLw0pkQm7sDS%2BR0%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2001%3A00%3A58%20GMT&X-Amz-Signature=58e5b73919ccbe4a0d29bb07a9bb2d60ddae56fcca303ba0ca4d0260a1c95f04'
$ws.Range("H3").Value = 20.627933959
$ws.Range("I3").Value = -103.29890938
$ws.Range("J3").Value = 0.7770282030105591

# Row 4
$ws.Range("A4").Value = '281474991265672-1749675495564'
$ws.Range("B4").Value = 'Harsh Brake'
$ws.Range("C4").Value = '2025-06-11T14:58:15.564'
$ws.Range("D4").Value = '''281474991265672'
$ws.Range("E4").Value = '''116'
$ws.Range("F4").Value = '''52215735'
$ws.Range("G4").Value = 'KEVIN DE LA O'
$ws.Range("K4").Value = 'No video URL'
$ws.Range("L4").Value = 'No video URL'
$ws.Range("H4").Value = 20.6922545
$ws.Range("I4").Value = -103.37642088
$ws.Range("J4").Value = 0.7694903612136841

# Row 5
$ws.Range("A5").Value = '281474991395097-1749663419318'
$ws.Range("B5").Value = 'Harsh Brake'
$ws.Range("C5").Value = '2025-06-11T11:36:59.318'
$ws.Range("D5").Value = '''281474991395097'
$ws.Range("E5").Value = '''125'
$ws.Range("F5").Value = '''51834055'
$ws.Range("G5").Value = 'DAVID SERRANO'
$ws.Range("K5").Value = 'No video URL'
$ws.Range("L5").Value = 'No video URL'
$ws.Range("H5").Value = 20.67376958
$ws.Range("I5").Value = -103.39834185
$ws.Range("J5").Value = 0.9524186849594116
